$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'66.584.29"
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.12%  '
$c = $ws.Range('D3')
$c.Value = "'3.522.93"
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.Value = "'607.24"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$c = $ws.Range('D6')
$c.Value = "'143.74"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.98%  '
$c = $ws.Range('D7')
$c.Value = "'3.522.35"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -3.01%  '
$c = $ws.Range('D8')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.22%  '
$c = $ws.Range('D9')
$c.Value = "'0.508"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +3.80%  '
$c = $ws.Range('D10')
$c.Value = "'7.70"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.64%  '
$ws.Range('E11').Value = '  -4.74%  '
$c = $ws.Range('D12')
$c.Value = "'0.406"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -2.90%  '
$c = $ws.Range('D13')
$c.Value = "'4.115.30"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.22%  '
$ws.Range('E14').Value = '  -6.76%  '
$c = $ws.Range('D15')
$c.Value = "'28.70"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -4.41%  '
$c = $ws.Range('D16')
$c.Value = "'3.526.37"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D17')
$c.Value = "'0.117"
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range('D18')
$c.Value = "'66.452.93"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.48%  '
$c = $ws.Range('D19')
$c.Value = "'10.79"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -7.19%  '
$c = $ws.Range('D20')
$c.Value = "'6.13"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -4.16%  '
$ws.Range('E21').Value = '  -3.50%  '
$c = $ws.Range('D22')
$c.Value = "'423.14"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.58%  '
$c = $ws.Range('D23')
$c.Value = "'0.590"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -5.18%  '
$c = $ws.Range('D24')
$c.Value = "'77.05"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.19%  '
$c = $ws.Range('D25')
$c.Value = "'3.669.41"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -6.81%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D28')
$c.Value = "'2.46"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D29')
$c.Value = "'7.86"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -6.36%  '
$c = $ws.Range('D30')
$c.Value = "'8.91"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -6.67%  '
$c = $ws.Range('D31')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.04%  '
$c = $ws.Range('D32')
$c.Value = "'3.529.25"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.89%  '
$c = $ws.Range('D33')
$c.Value = "'0.154"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.08%  '
$c = $ws.Range('D34')
$c.Value = "'24.21"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -5.21%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -9.98%  '
$c = $ws.Range('D37')
$c.Value = "'7.55"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -4.51%  '
$ws.Range('E38').Value = '  -5.17%  '
$c = $ws.Range('D39')
$c.Value = "'173.53"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$c = $ws.Range('D40')
$c.Value = "'5.20"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -8.73%  '
$c = $ws.Range('D41')
$c.Value = "'0.0811"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -6.04%  '
$ws.Range('E42').Value = '  -5.43%  '
$ws.Range('E43').Value = '  -5.56%  '
$c = $ws.Range('D44')
$c.Value = "'45.49"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('E45').Value = '  -6.89%  '
$ws.Range('E46').Value = '  +0.03%  '
$c = $ws.Range('D47')
$c.Value = "'2.37"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -8.71%  '
$ws.Range('E48').Value = '  -2.15%  '
$c = $ws.Range('D49')
$c.Value = "'1.11"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -5.77%  '
$c = $ws.Range('D50')
$c.Value = "'22.91"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -4.27%  '
$ws.Range('E51').Value = '  -5.93%  '
